$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.463.37"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.917.62"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.67%  "
$ws.Range("D5").Value = "'325.53"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'0.4810"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "'0.4051"
$ws.Range("E8").Value = "  -0.58%  "
$ws.Range("D9").Value = "'0.08208"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").Value = "'23.32"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "1.907.45"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'6.051"
$ws.Range("D14").Value = "'7.224"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "'91.19"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'0.06839"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'0.00001038"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'17.51"
$ws.Range("E19").Value = "  -1.02%  "
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("D21").Value = "29.475.28"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "'5.660"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").Value = "'11.84"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +1.47%  "
$ws.Range("D25").Value = "2.133.10"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'6.614"
$ws.Range("E26").Value = "  +6.12%  "
$ws.Range("D27").Value = "'155.81"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("D28").Value = "'19.99"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'120.29"
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("D31").Value = "'1.011"
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").Value = "'0.09605"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").Value = "'3.556"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "'1.370"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "'0.06261"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("D37").Value = "'0.02282"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'1.180"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "'0.5924"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("D40").Value = "'10.72"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'7.856"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'0.1844"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'2.387"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "'12.43"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'1.930"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("D50").Value = "'117.93"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "'2.419"
$ws.Range("E51").Value = "  +3.13%  "
